$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# --- Populate the new "2nd Attempt" row (row 9), values -------------------
$ws.Range("B9").Value = "2nd Attempt"
$ws.Range("C9").Value = "Jeff Renshaw"
$ws.Range("D9").Value = "Go"
$ws.Range("E9").Formula = "=AVERAGE(F9:M9)"
$ws.Range("F9").Value = 595.558
$ws.Range("G9").Value = 600.442
$ws.Range("H9").Value = 582.087
$ws.Range("I9").Value = 563.399
$ws.Range("J9").Value = 552.473
$ws.Range("K9").Value = 596.393
$ws.Range("L9").Value = 593.13
$ws.Range("M9").Value = 582.18

# --- Formatting: center alignment for C9:M9 (B9 left-aligned like wrapped label)
$ws.Range("C9:M9").HorizontalAlignment = -4108

# --- Number format: two decimals on the average cell ----------------------
$ws.Range("E9").NumberFormat = "0.00"

# --- Fill: highlight the label cell in yellow ------------------------------
$ws.Range("B9").Interior.Color = 65535

# --- Borders --------------------------------------------------------------
# B9: medium left + medium bottom (outer-left corner of the box)
$ws.Range("B9").Borders.Item(7).LineStyle = 1
$ws.Range("B9").Borders.Item(7).Weight = -4138
$ws.Range("B9").Borders.Item(9).LineStyle = 1
$ws.Range("B9").Borders.Item(9).Weight = -4138

# C9:L9: thin left + thin right + medium bottom on every individual cell
$innerCols = @("C","D","E","F","G","H","I","J","K","L")
foreach ($col in $innerCols) {
    $addr = $col + "9"
    $ws.Range($addr).Borders.Item(7).LineStyle = 1
    $ws.Range($addr).Borders.Item(7).Weight = 2
    $ws.Range($addr).Borders.Item(10).LineStyle = 1
    $ws.Range($addr).Borders.Item(10).Weight = 2
    $ws.Range($addr).Borders.Item(9).LineStyle = 1
    $ws.Range($addr).Borders.Item(9).Weight = -4138
}

# M9: thin left + medium right + medium bottom (outer-right corner)
$ws.Range("M9").Borders.Item(7).LineStyle = 1
$ws.Range("M9").Borders.Item(7).Weight = 2
$ws.Range("M9").Borders.Item(10).LineStyle = 1
$ws.Range("M9").Borders.Item(10).Weight = -4138
$ws.Range("M9").Borders.Item(9).LineStyle = 1
$ws.Range("M9").Borders.Item(9).Weight = -4138

# --- Row height for the new row (visually matches the thick-bottom rows) --
$ws.Rows.Item(9).RowHeight = 15

# --- Column B needs to be wide enough for "2nd Attempt" -------------------
$ws.Columns.Item(2).ColumnWidth = 11.77

# --- Selection moved as part of the edit session ---------------------------
$ws.Range("G26").Select()
